$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 43).Value = 0.71403170592694076
$ws.Cells.Item(1, 64).Value = 0.97720942995708926
$ws.Cells.Item(2, 4).Value = 0.97651111671114255
$ws.Cells.Item(3, 1).Value = 0.99380244616658719
$ws.Cells.Item(3, 2).Value = 0.91898594927626343
$ws.Cells.Item(4, 3).Value = 0.97159679339031357
$ws.Cells.Item(4, 66).Value = 0.83095857089805469
$ws.Cells.Item(5, 3).Value = 0.99689713480146014
$ws.Cells.Item(5, 4).Value = 0.9499643794600221
$ws.Cells.Item(5, 6).Value = 0.79710640515974518
$ws.Cells.Item(5, 7).Value = 0.97918956265191048
$ws.Cells.Item(6, 28).Value = 0.7716412865800566
$ws.Cells.Item(7, 9).Value = 0.94134327713637567
$ws.Cells.Item(7, 64).Value = 0.7657603063751135
$ws.Cells.Item(7, 68).Value = 0.5384347030897858
$ws.Cells.Item(8, 6).Value = 0.98158649439654999
$ws.Cells.Item(9, 8).Value = 0.96272324530105235
$ws.Cells.Item(9, 10).Value = 0.89800964226414459
$ws.Cells.Item(10, 8).Value = 0.63294864873697443
$ws.Cells.Item(10, 11).Value = 0.91495215406785335
$ws.Cells.Item(10, 12).Value = 0.96104681569718831
$ws.Cells.Item(11, 9).Value = 0.94630246314014732
$ws.Cells.Item(11, 13).Value = 0.85901920132283882
$ws.Cells.Item(11, 31).Value = 0.94229862725191671
$ws.Cells.Item(11, 57).Value = 0.87145845806477851
$ws.Cells.Item(12, 14).Value = 0.86296876152739588
$ws.Cells.Item(13, 12).Value = 0.88169052953296534
$ws.Cells.Item(13, 15).Value = 0.82814304942099048
$ws.Cells.Item(13, 25).Value = 0.93546283476005465
$ws.Cells.Item(14, 8).Value = 0.97726033249680166
$ws.Cells.Item(14, 29).Value = 0.78476570497455356
$ws.Cells.Item(15, 34).Value = 0.89931760816474848
$ws.Cells.Item(16, 6).Value = 0.64989154722941667
$ws.Cells.Item(16, 14).Value = 0.74402968981339856
$ws.Cells.Item(16, 15).Value = 0.79369620009769182
$ws.Cells.Item(17, 16).Value = 0.82213349621560083
$ws.Cells.Item(17, 18).Value = 0.99656010715271404
$ws.Cells.Item(17, 38).Value = 0.79801151904710532
$ws.Cells.Item(18, 20).Value = 0.8187728058227004
$ws.Cells.Item(19, 18).Value = 0.80105997657266492
$ws.Cells.Item(20, 19).Value = 0.97492634520640165
$ws.Cells.Item(20, 21).Value = 0.8010201253675322
$ws.Cells.Item(21, 19).Value = 0.97362429680847473
$ws.Cells.Item(21, 32).Value = 0.98059398845780632
$ws.Cells.Item(21, 45).Value = 0.92209655288556469
$ws.Cells.Item(22, 20).Value = 0.76494206892772287
$ws.Cells.Item(22, 21).Value = 0.88369903975777819
$ws.Cells.Item(22, 24).Value = 0.82024233591779039
$ws.Cells.Item(23, 22).Value = 0.94639978594477592
$ws.Cells.Item(24, 23).Value = 0.91713000562341063
$ws.Cells.Item(24, 59).Value = 0.69564465455974966
$ws.Cells.Item(25, 23).Value = 0.82001590668299917
$ws.Cells.Item(26, 25).Value = 0.97757115744987533
$ws.Cells.Item(27, 25).Value = 0.87869648136617817
$ws.Cells.Item(27, 26).Value = 0.92054388176400692
$ws.Cells.Item(28, 26).Value = 0.809336152085259
$ws.Cells.Item(28, 27).Value = 0.61539579162460534
$ws.Cells.Item(29, 27).Value = 0.93072013929115671
$ws.Cells.Item(29, 28).Value = 0.99468173845968122
$ws.Cells.Item(30, 28).Value = 0.98940268609875059
$ws.Cells.Item(30, 31).Value = 0.55482914298375763
$ws.Cells.Item(30, 32).Value = 0.8177024799586301
$ws.Cells.Item(32, 26).Value = 0.86565459448331172
$ws.Cells.Item(32, 58).Value = 0.93353494721574592
$ws.Cells.Item(33, 31).Value = 0.90619168977617837
$ws.Cells.Item(33, 35).Value = 0.70357066747376207
$ws.Cells.Item(34, 32).Value = 0.97121115889974852
$ws.Cells.Item(34, 33).Value = 0.89750841222734257
$ws.Cells.Item(36, 29).Value = 0.96544774658021881
$ws.Cells.Item(36, 34).Value = 0.92231762534292283
$ws.Cells.Item(36, 35).Value = 0.61731273235003936
$ws.Cells.Item(37, 39).Value = 0.86414242385552975
$ws.Cells.Item(38, 29).Value = 0.97055840625492129
$ws.Cells.Item(38, 36).Value = 0.82748773466445424
$ws.Cells.Item(38, 37).Value = 0.93628613457163667
$ws.Cells.Item(39, 38).Value = 0.91250236288909969
$ws.Cells.Item(39, 49).Value = 0.91340391967738133
$ws.Cells.Item(40, 25).Value = 0.9468513923818529
$ws.Cells.Item(40, 27).Value = 0.92175908770683612
$ws.Cells.Item(40, 38).Value = 0.62840475633037207
$ws.Cells.Item(40, 39).Value = 0.92684240227354775
$ws.Cells.Item(41, 39).Value = 0.98615671880559419
$ws.Cells.Item(41, 42).Value = 0.90076068442186807
$ws.Cells.Item(41, 43).Value = 0.97428434997633029
$ws.Cells.Item(43, 35).Value = 0.93675338028019617
$ws.Cells.Item(43, 42).Value = 0.66286360016570955
$ws.Cells.Item(43, 45).Value = 0.87240994951815143
$ws.Cells.Item(44, 42).Value = 0.7264733862682875
$ws.Cells.Item(44, 43).Value = 0.85438752556072373
$ws.Cells.Item(44, 45).Value = 0.58421234517551501
$ws.Cells.Item(44, 46).Value = 0.85683983638162908
$ws.Cells.Item(45, 51).Value = 0.99841892731765824
$ws.Cells.Item(46, 45).Value = 0.7670561656493512
$ws.Cells.Item(46, 48).Value = 0.65321689236016423
$ws.Cells.Item(47, 26).Value = 0.73444551267948377
$ws.Cells.Item(47, 46).Value = 0.84336875616624818
$ws.Cells.Item(47, 49).Value = 0.92791653527144102
$ws.Cells.Item(48, 7).Value = 0.99287935847118991
$ws.Cells.Item(48, 24).Value = 0.67479280735778957
$ws.Cells.Item(49, 51).Value = 0.98131681620891342
$ws.Cells.Item(50, 48).Value = 0.74642678876582347
$ws.Cells.Item(50, 51).Value = 0.63166854372836334
$ws.Cells.Item(52, 45).Value = 0.98866599095255769
$ws.Cells.Item(52, 50).Value = 0.96311554673882016
$ws.Cells.Item(53, 51).Value = 0.87604536631793684
$ws.Cells.Item(53, 52).Value = 0.97508311771472977
$ws.Cells.Item(53, 65).Value = 0.92657937928518708
$ws.Cells.Item(54, 55).Value = 0.84899549577922795
$ws.Cells.Item(55, 56).Value = 0.90734968146235473
$ws.Cells.Item(55, 57).Value = 0.91074741798569447
$ws.Cells.Item(56, 54).Value = 0.83633385883594036
$ws.Cells.Item(56, 57).Value = 0.91129212539136262
$ws.Cells.Item(56, 58).Value = 0.96571211321581374
$ws.Cells.Item(57, 62).Value = 0.9258377414690534
$ws.Cells.Item(58, 60).Value = 0.79073977968261799
$ws.Cells.Item(59, 53).Value = 0.99809032821890087
$ws.Cells.Item(59, 61).Value = 0.67249957114915637
$ws.Cells.Item(60, 59).Value = 0.98526481623615703
$ws.Cells.Item(60, 61).Value = 0.67100643651922942
$ws.Cells.Item(61, 62).Value = 0.51241644686344179
$ws.Cells.Item(62, 12).Value = 0.68114749142279551
$ws.Cells.Item(62, 34).Value = 0.92752872197023539
$ws.Cells.Item(62, 60).Value = 0.96582138677041418
$ws.Cells.Item(62, 64).Value = 0.85853818189086528
$ws.Cells.Item(63, 61).Value = 0.60121123649746733
$ws.Cells.Item(63, 65).Value = 0.93897052273016723
$ws.Cells.Item(64, 63).Value = 0.85706621768748259
$ws.Cells.Item(64, 66).Value = 0.97958943298950252
$ws.Cells.Item(65, 31).Value = 0.89832509211338518
$ws.Cells.Item(65, 66).Value = 0.99313793679051532
$ws.Cells.Item(66, 68).Value = 0.63069019655098946
$ws.Cells.Item(67, 1).Value = 0.52975439074784181
$ws.Cells.Item(67, 65).Value = 0.97204627321028858
$ws.Cells.Item(67, 66).Value = 0.69384978768741856
$ws.Cells.Item(67, 68).Value = 0.77598657006122185
$ws.Cells.Item(68, 55).Value = 0.8417949359619461

Write-Host "Applied 136 cell updates"
